$wb = $excel.ActiveWorkbook

# Sheet 1 = "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(3, 7).Value = 99

$ws1.Cells.Item(6, 6).Value = 2997
$ws1.Cells.Item(7, 6).Value = 1691
$ws1.Cells.Item(8, 6).Value = 2015

$ws1.Cells.Item(9, 6).Value = 320
$ws1.Cells.Item(9, 7).Value = 218

$ws1.Cells.Item(11, 6).Value = 865
$ws1.Cells.Item(12, 6).Value = 953
$ws1.Cells.Item(13, 6).Value = 203
$ws1.Cells.Item(14, 6).Value = 421

$ws1.Cells.Item(18, 6).Value = 534
$ws1.Cells.Item(19, 6).Value = 7261

$ws1.Cells.Item(21, 6).Value = 1983
$ws1.Cells.Item(22, 6).Value = 190

$ws1.Cells.Item(25, 6).Value = 448
$ws1.Cells.Item(26, 6).Value = 501

$ws1.Cells.Item(28, 6).Value = 1121
$ws1.Cells.Item(29, 6).Value = 946

$ws1.Cells.Item(31, 7).Value = 65

$ws1.Cells.Item(33, 6).Value = 1122

$ws1.Cells.Item(36, 6).Value = 14

$ws1.Cells.Item(38, 6).Value = 255

$ws1.Cells.Item(40, 6).Value = 151
$ws1.Cells.Item(41, 6).Value = 279

# Sheet 4 = "全部类型" (All Types) - mirrors the same events
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(3, 7).Value = 99

$ws4.Cells.Item(9, 6).Value = 2997
$ws4.Cells.Item(10, 6).Value = 1691
$ws4.Cells.Item(11, 6).Value = 2015

$ws4.Cells.Item(12, 6).Value = 320
$ws4.Cells.Item(12, 7).Value = 218

$ws4.Cells.Item(14, 6).Value = 865
$ws4.Cells.Item(16, 6).Value = 953
$ws4.Cells.Item(17, 6).Value = 203
$ws4.Cells.Item(18, 6).Value = 421

$ws4.Cells.Item(22, 6).Value = 534
$ws4.Cells.Item(23, 6).Value = 7261

$ws4.Cells.Item(25, 6).Value = 1983
$ws4.Cells.Item(27, 6).Value = 190

$ws4.Cells.Item(30, 6).Value = 448
$ws4.Cells.Item(31, 6).Value = 501

$ws4.Cells.Item(33, 6).Value = 1121
$ws4.Cells.Item(34, 6).Value = 946

$ws4.Cells.Item(36, 7).Value = 65

$ws4.Cells.Item(37, 6).Value = 1122

$ws4.Cells.Item(40, 6).Value = 14

$ws4.Cells.Item(42, 6).Value = 255

$ws4.Cells.Item(44, 6).Value = 151
$ws4.Cells.Item(45, 6).Value = 279
